$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 33, item id 5512
$ws.Range("H33").Value = 525.6667
$ws.Range("I33").Value = 525.6667
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 525.6667
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -296.6667

$ws = $wb.Worksheets.Item("ALC")  # row 39, item id 4603
$ws.Range("H39").Value = 2187.3845
$ws.Range("I39").Value = 2486.5557
$ws.Range("J39").Value = 1514.25
$ws.Range("K39").Value = 7459.6671
$ws.Range("L39").Value = 4542.75
$ws.Range("M39").Value = -7163.6671
$ws.Range("N39").Value = -5134.75

$ws = $wb.Worksheets.Item("ALC")  # row 41, item id 5478
$ws.Range("H41").Value = 487.5
$ws.Range("I41").Value = 487.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 487.5
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -47.5

$ws = $wb.Worksheets.Item("ALC")  # row 92, item id 19901
$ws.Range("H92").Value = 1169.7894
$ws.Range("I92").Value = 1067.25
$ws.Range("J92").Value = 1716.6666
$ws.Range("K92").Value = 1067.25
$ws.Range("L92").Value = 1716.6666
$ws.Range("M92").Value = 180.75
$ws.Range("N92").Value = -4212.6666

$ws = $wb.Worksheets.Item("ALC")  # row 94, item id 19905
$ws.Range("H94").Value = 2699.5
$ws.Range("I94").Value = 2699.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2699.5
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2248.5

$ws = $wb.Worksheets.Item("ALC")  # row 100, item id 19906
$ws.Range("H100").Value = 1133
$ws.Range("I100").Value = 1041.25
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1041.25
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -500.25
$ws.Range("N100").Value = -2582

$ws = $wb.Worksheets.Item("ALC")  # row 106, item id 19903
$ws.Range("H106").Value = 8664
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 8664
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 8664
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -9926

$ws = $wb.Worksheets.Item("ALC")  # row 132, item id 44049
$ws.Range("H132").Value = 1081601
$ws.Range("I132").Value = 2149936.2
$ws.Range("J132").Value = 13265.833
$ws.Range("K132").Value = 6449808.600000001
$ws.Range("L132").Value = 39797.499
$ws.Range("M132").Value = -6447278.600000001
$ws.Range("N132").Value = -44857.499

$ws = $wb.Worksheets.Item("ALC")  # row 135, item id 44047
$ws.Range("H135").Value = 4032.2163
$ws.Range("I135").Value = 1165.08
$ws.Range("J135").Value = 10005.417
$ws.Range("K135").Value = 10485.72
$ws.Range("L135").Value = 90048.753
$ws.Range("M135").Value = -7950.719999999999
$ws.Range("N135").Value = -95118.753

$ws = $wb.Worksheets.Item("ALC")  # row 141, item id 44161
$ws.Range("H141").Value = 2884.0789
$ws.Range("I141").Value = 2708.2856
$ws.Range("J141").Value = 4935
$ws.Range("K141").Value = 8124.8568
$ws.Range("L141").Value = 14805
$ws.Range("M141").Value = -2944.8568
$ws.Range("N141").Value = -25165

$ws = $wb.Worksheets.Item("ARM")  # row 32, item id 44147
$ws.Range("H32").Value = 1712294.5
$ws.Range("I32").Value = 1139.8959
$ws.Range("J32").Value = 6543790
$ws.Range("K32").Value = 1139.8959
$ws.Range("L32").Value = 6543790
$ws.Range("M32").Value = -852.8959
$ws.Range("N32").Value = -6544364

$ws = $wb.Worksheets.Item("ARM")  # row 63, item id 12528
$ws.Range("H63").Value = 6620
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 6620
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 6620
$ws.Range("N63").Value = -7992

$ws = $wb.Worksheets.Item("ARM")  # row 66, item id 12528
$ws.Range("H66").Value = 6620
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 6620
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 33100
$ws.Range("N66").Value = -39964

$ws = $wb.Worksheets.Item("ARM")  # row 97, item id 19941
$ws.Range("H97").Value = 1991
$ws.Range("I97").Value = 1624.5
$ws.Range("J97").Value = 2409.8572
$ws.Range("K97").Value = 1624.5
$ws.Range("L97").Value = 2409.8572
$ws.Range("M97").Value = -1128.5
$ws.Range("N97").Value = -3401.8572

$ws = $wb.Worksheets.Item("ARM")  # row 132, item id 43997
$ws.Range("H132").Value = 748645.9
$ws.Range("I132").Value = 892720.9
$ws.Range("J132").Value = 136327.12
$ws.Range("K132").Value = 2678162.7
$ws.Range("L132").Value = 408981.36
$ws.Range("M132").Value = -2675632.7
$ws.Range("N132").Value = -414041.36

$ws = $wb.Worksheets.Item("BSM")  # row 86, item id 12526
$ws.Range("H86").Value = 15849.667
$ws.Range("I86").Value = 35333.332
$ws.Range("J86").Value = 6107.8335
$ws.Range("K86").Value = 35333.332
$ws.Range("L86").Value = 6107.8335
$ws.Range("M86").Value = -34210.332
$ws.Range("N86").Value = -8353.833500000001

$ws = $wb.Worksheets.Item("BSM")  # row 89, item id 12526
$ws.Range("H89").Value = 15849.667
$ws.Range("I89").Value = 35333.332
$ws.Range("J89").Value = 6107.8335
$ws.Range("K89").Value = 176666.66
$ws.Range("L89").Value = 30539.1675
$ws.Range("M89").Value = -171050.66
$ws.Range("N89").Value = -41771.1675

$ws = $wb.Worksheets.Item("BSM")  # row 94, item id 19939
$ws.Range("H94").Value = 5209.1333
$ws.Range("I94").Value = 2626.5806
$ws.Range("J94").Value = 10927.643
$ws.Range("K94").Value = 2626.5806
$ws.Range("L94").Value = 10927.643
$ws.Range("M94").Value = -2175.5806
$ws.Range("N94").Value = -11829.643

$ws = $wb.Worksheets.Item("CRP")  # row 22, item id 5367
$ws.Range("H22").Value = 1090504.9
$ws.Range("I22").Value = 1716148.9
$ws.Range("J22").Value = 107349.86
$ws.Range("K22").Value = 1716148.9
$ws.Range("L22").Value = 107349.86
$ws.Range("M22").Value = -1715798.9
$ws.Range("N22").Value = -108049.86

$ws = $wb.Worksheets.Item("CRP")  # row 31, item id 44023
$ws.Range("H31").Value = 4015.3794
$ws.Range("I31").Value = 1878.1052
$ws.Range("J31").Value = 8076.2
$ws.Range("K31").Value = 1878.1052
$ws.Range("L31").Value = 8076.2
$ws.Range("M31").Value = -1583.1052
$ws.Range("N31").Value = -8666.200000000001

$ws = $wb.Worksheets.Item("CRP")  # row 34, item id 44023
$ws.Range("H34").Value = 4015.3794
$ws.Range("I34").Value = 1878.1052
$ws.Range("J34").Value = 8076.2
$ws.Range("K34").Value = 1878.1052
$ws.Range("L34").Value = 8076.2
$ws.Range("M34").Value = -1676.1052
$ws.Range("N34").Value = -8480.200000000001

$ws = $wb.Worksheets.Item("CRP")  # row 62, item id 12580
$ws.Range("H62").Value = 7898.3335
$ws.Range("I62").Value = 16923.5
$ws.Range("J62").Value = 3385.75
$ws.Range("K62").Value = 16923.5
$ws.Range("L62").Value = 3385.75
$ws.Range("M62").Value = -16299.5
$ws.Range("N62").Value = -4633.75

$ws = $wb.Worksheets.Item("CRP")  # row 65, item id 12580
$ws.Range("H65").Value = 7898.3335
$ws.Range("I65").Value = 16923.5
$ws.Range("J65").Value = 3385.75
$ws.Range("K65").Value = 84617.5
$ws.Range("L65").Value = 16928.75
$ws.Range("M65").Value = -81497.5
$ws.Range("N65").Value = -23168.75

$ws = $wb.Worksheets.Item("CRP")  # row 68, item id 10611
$ws.Range("H68").Value = 67498.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 67498.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 67498.5
$ws.Range("N68").Value = -68996.5

$ws = $wb.Worksheets.Item("CRP")  # row 71, item id 10611
$ws.Range("H71").Value = 67498.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 67498.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 202495.5
$ws.Range("N71").Value = -209983.5

$ws = $wb.Worksheets.Item("CRP")  # row 132, item id 44019
$ws.Range("H132").Value = 5612.0303
$ws.Range("I132").Value = 4275.8276
$ws.Range("J132").Value = 15299.5
$ws.Range("K132").Value = 12827.4828
$ws.Range("L132").Value = 45898.5
$ws.Range("M132").Value = -10297.4828
$ws.Range("N132").Value = -50958.5

$ws = $wb.Worksheets.Item("CUL")  # row 14, item id 12886
$ws.Range("H14").Value = 299.64706
$ws.Range("I14").Value = 299.64706
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 898.94118
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -725.94118

$ws = $wb.Worksheets.Item("CUL")  # row 52, item id 31902
$ws.Range("H52").Value = 2315518.2
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 2315518.2
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 6946554.600000001
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -6947086.600000001

$ws = $wb.Worksheets.Item("CUL")  # row 109, item id 27854
$ws.Range("H109").Value = 12600.667
$ws.Range("I109").Value = 9445.556
$ws.Range("J109").Value = 17333.334
$ws.Range("K109").Value = 28336.668
$ws.Range("L109").Value = 52000.00199999999
$ws.Range("M109").Value = -27296.668
$ws.Range("N109").Value = -54080.00199999999

$ws = $wb.Worksheets.Item("CUL")  # row 132, item id 43972
$ws.Range("H132").Value = 2248.15
$ws.Range("I132").Value = 1042.8334
$ws.Range("J132").Value = 4056.125
$ws.Range("K132").Value = 9385.500599999999
$ws.Range("L132").Value = 36505.125
$ws.Range("M132").Value = -6855.500599999999
$ws.Range("N132").Value = -41565.125

$ws = $wb.Worksheets.Item("CUL")  # row 136, item id 44093
$ws.Range("H136").Value = 16670101
$ws.Range("I136").Value = 23812488
$ws.Range("J136").Value = 4531
$ws.Range("K136").Value = 71437464
$ws.Range("L136").Value = 13593
$ws.Range("M136").Value = -71432364
$ws.Range("N136").Value = -23793

$ws = $wb.Worksheets.Item("GSM")  # row 132, item id 44008
$ws.Range("H132").Value = 66670224
$ws.Range("I132").Value = 100003060
$ws.Range("J132").Value = 4560.2
$ws.Range("K132").Value = 300009180
$ws.Range("L132").Value = 13680.6
$ws.Range("M132").Value = -300006650
$ws.Range("N132").Value = -18740.6

$ws = $wb.Worksheets.Item("LTW")  # row 68, item id 12563
$ws.Range("H68").Value = 2155.037
$ws.Range("I68").Value = 1914.5
$ws.Range("J68").Value = 2842.2856
$ws.Range("K68").Value = 1914.5
$ws.Range("L68").Value = 2842.2856
$ws.Range("M68").Value = -1165.5
$ws.Range("N68").Value = -4340.2856

$ws = $wb.Worksheets.Item("LTW")  # row 71, item id 12563
$ws.Range("H71").Value = 2155.037
$ws.Range("I71").Value = 1914.5
$ws.Range("J71").Value = 2842.2856
$ws.Range("K71").Value = 9572.5
$ws.Range("L71").Value = 14211.428
$ws.Range("M71").Value = -5828.5
$ws.Range("N71").Value = -21699.428

$ws = $wb.Worksheets.Item("LTW")  # row 106, item id 18713
$ws.Range("H106").Value = 28842.143
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 28842.143
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 28842.143
$ws.Range("N106").Value = -31366.143

$ws = $wb.Worksheets.Item("LTW")  # row 132, item id 44058
$ws.Range("H132").Value = 3256.0322
$ws.Range("I132").Value = 3229.4644
$ws.Range("J132").Value = 3504
$ws.Range("K132").Value = 9688.393199999999
$ws.Range("L132").Value = 10512
$ws.Range("M132").Value = -7158.393199999999
$ws.Range("N132").Value = -15572

$ws = $wb.Worksheets.Item("WVR")  # row 62, item id 12589
$ws.Range("H62").Value = 12799.777
$ws.Range("I62").Value = 15417.333
$ws.Range("J62").Value = 11491
$ws.Range("K62").Value = 15417.333
$ws.Range("L62").Value = 11491
$ws.Range("M62").Value = -14793.333
$ws.Range("N62").Value = -12739

$ws = $wb.Worksheets.Item("WVR")  # row 65, item id 12589
$ws.Range("H65").Value = 12799.777
$ws.Range("I65").Value = 15417.333
$ws.Range("J65").Value = 11491
$ws.Range("K65").Value = 77086.66500000001
$ws.Range("L65").Value = 57455
$ws.Range("M65").Value = -73966.66500000001
$ws.Range("N65").Value = -63695

$ws = $wb.Worksheets.Item("WVR")  # row 100, item id 19981
$ws.Range("H100").Value = 1796.8572
$ws.Range("I100").Value = 816.1667
$ws.Range("J100").Value = 2532.375
$ws.Range("K100").Value = 1632.3334
$ws.Range("L100").Value = 5064.75
$ws.Range("M100").Value = -1091.3334
$ws.Range("N100").Value = -6146.75

$ws = $wb.Worksheets.Item("WVR")  # row 104, item id 18691
$ws.Range("H104").Value = 49746
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 49746
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 49746
$ws.Range("N104").Value = -56734

$ws = $wb.Worksheets.Item("WVR")  # row 132, item id 44029
$ws.Range("H132").Value = 5917.5273
$ws.Range("I132").Value = 4158.915
$ws.Range("J132").Value = 16249.375
$ws.Range("K132").Value = 12476.745
$ws.Range("L132").Value = 48748.125
$ws.Range("M132").Value = -9946.744999999999
$ws.Range("N132").Value = -53808.125
